$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: the three data rows get reordered (the 83b4f57b row moves
# from row 2 down to row 4) and the 83b4f57b row's status/date is refreshed
# to reflect that it is now ready for a new handoff.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "ffffefd082a9-af3c-47af-88e1-af967a44b92d.md"
$ws.Range("B2").Value = "e2e\ffffefd082a9-af3c-47af-88e1-af967a44b92d.md"
$ws.Range("E2").Value = "Handed back: in sync with en-US"
$ws.Range("F2").Value = "Handed back: in sync with en-US"
$ws.Range("G2").Value = "2016-08-20 13:03:55"

$ws.Range("A3").Value = "ffffffce63173e-171c-44a5-8f2d-3ac06ee6e219.md"
$ws.Range("B3").Value = "e2e\ffffffce63173e-171c-44a5-8f2d-3ac06ee6e219.md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-08-20 13:03:55"

$ws.Range("A4").Value = "83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.md"
$ws.Range("B4").Value = "e2e\83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.md"
$ws.Range("E4").Value = "Ready for handoff"
$ws.Range("F4").Value = "Ready for handoff"
$ws.Range("G4").Value = "2016-08-20 13:05:44"

# ---------------------------------------------------------------------------
# zh-cn sheet: same row reshuffle, localized to the zh-cn xliff filenames.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "ffffefd082a9-af3c-47af-88e1-af967a44b92d.md"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("F2").Value = "False"
$wsZh.Range("G2").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.074d7c1959bb106be22360d7b6cb090df0a970a1.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-20 13:03:51"
$wsZh.Range("I2").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.md"
$wsZh.Range("J2").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.074d7c1959bb106be22360d7b6cb090df0a970a1.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-20 13:04:12"
$wsZh.Range("P2").Value = ""

$wsZh.Range("A3").Value = "ffffffce63173e-171c-44a5-8f2d-3ac06ee6e219.md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.074d7c1959bb106be22360d7b6cb090df0a970a1.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-20 13:03:51"
$wsZh.Range("I3").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.md"
$wsZh.Range("J3").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.074d7c1959bb106be22360d7b6cb090df0a970a1.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-20 13:04:12"
$wsZh.Range("P3").Value = ""

$wsZh.Range("A4").Value = "83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = "83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.e9100736e12a57bd6ecb02489df9499d164b78a3.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-08-20 13:05:40"
$wsZh.Range("I4").Value = "83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.md"
$wsZh.Range("J4").Value = "83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.e9100736e12a57bd6ecb02489df9499d164b78a3.zh-cn.xlf"
$wsZh.Range("K4").Value = "2016-08-20 13:05:15"
$wsZh.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58ce6d650190a178daa803a63a606a40eb95cb29/e2e/83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d902cea256ba90afe5675decb2fa335d9bbb7cc/e2e/83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.md."

$wsZh.Columns.Item(16).ColumnWidth = 40

# ---------------------------------------------------------------------------
# de-de sheet: same row reshuffle, localized to the de-de xliff filenames.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "ffffefd082a9-af3c-47af-88e1-af967a44b92d.md"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("F2").Value = "False"
$wsDe.Range("G2").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.074d7c1959bb106be22360d7b6cb090df0a970a1.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-20 13:03:55"
$wsDe.Range("I2").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.md"
$wsDe.Range("J2").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.074d7c1959bb106be22360d7b6cb090df0a970a1.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-20 13:04:19"
$wsDe.Range("P2").Value = ""

$wsDe.Range("A3").Value = "ffffffce63173e-171c-44a5-8f2d-3ac06ee6e219.md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.074d7c1959bb106be22360d7b6cb090df0a970a1.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-20 13:03:55"
$wsDe.Range("I3").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.md"
$wsDe.Range("J3").Value = "9a023f4c-dea3-4eab-91ea-79f72a209048.074d7c1959bb106be22360d7b6cb090df0a970a1.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-20 13:04:19"
$wsDe.Range("P3").Value = ""

$wsDe.Range("A4").Value = "83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = "83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.e9100736e12a57bd6ecb02489df9499d164b78a3.de-de.xlf"
$wsDe.Range("H4").Value = "2016-08-20 13:05:44"
$wsDe.Range("I4").Value = "83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.md"
$wsDe.Range("J4").Value = "83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.e9100736e12a57bd6ecb02489df9499d164b78a3.de-de.xlf"
$wsDe.Range("K4").Value = "2016-08-20 13:05:22"
$wsDe.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58ce6d650190a178daa803a63a606a40eb95cb29/e2e/83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d902cea256ba90afe5675decb2fa335d9bbb7cc/e2e/83b4f57b-b40d-4be4-a004-e0dfed0ccd4e.md."

$wsDe.Columns.Item(16).ColumnWidth = 40
